$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy rows 95:96 and insert the copy before row 95, shifting existing rows down.
$ws.Range("A95:T96").Copy() | Out-Null
$ws.Range("A95:T96").Insert() | Out-Null

# Update the new row 95 and row 96 with their new values.
$ws.Range("D95").Value = 44904
$ws.Range("M95").Value = 1000
$ws.Range("N95").Value = 4800
$ws.Range("O95").Value = 5000
$ws.Range("P95").Value = 4910
$ws.Range("S95").Value = 1228

$ws.Range("D96").Value = 44904
$ws.Range("M96").Value = 800
$ws.Range("N96").Value = 4800
$ws.Range("O96").Value = 5000
$ws.Range("P96").Value = 4888
$ws.Range("S96").Value = 1222
